# Apply the edits described by the diff to "Ver avaliações.xlsx"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix typo in D9: "dos todos" -> "de todos"
$ws.Range("D9").Value = "Apresenta todos as Avaliações de todos os seus serviços"

# Move the active selection from C9 to D9 (matches sheetView selection change)
$ws.Range("D9").Select()
